# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Update summary figures (Valor Mora total, worker/period counts) ---
$ws.Range("E11").Value = 284700
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 5

# --- Grab the formatting of the old table's closing row (23) before we
# start rewriting rows, so it can be re-applied to the new closing row. ---
$ws.Range("B23:J23").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Replace employee data table ---
# Rows 16-19 keep worker 8055429 / OMAR DAVID RODRIGUEZ VALENCIA but with
# periods 2504..2507 (ascending) and the new "Valor Mora" amount (G).
$periods = @("2504", "2505", "2506", "2507", "2508")
for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = 16 + $i
    $ws.Cells.Item($r, 2).Value = "CC"
    $ws.Cells.Item($r, 3).Value = "8055429"
    $ws.Cells.Item($r, 4).Value = "OMAR DAVID RODRIGUEZ VALENCIA"
    $ws.Cells.Item($r, 5).Value = $periods[$i]
    $ws.Cells.Item($r, 6).Value = 56940
    $ws.Cells.Item($r, 7).Value = 1423500
}

# Delete the now-redundant old rows 21-23 (leftover duplicate worker rows);
# this shifts the signature block (rows 28-29) up to rows 25-26.
$ws.Rows("21:23").Delete()
